$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 10137753.70137369
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("G2").Value = 208739759.1852868
